# Update the "target_data" sheet: fill M2:O29 with 200 and set the
# resulting selection to match (anchor M2, range M2:O29).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("target_data")
$ws.Activate()

$ws.Range("M2:O29").Value = 200

$ws.Range("M2:O29").Select()
